$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Updated daily-case counts (col B), revised figures from the 15-Apr CDC pull ---
$ws.Range("B27").Value = 6
$ws.Range("B30").Value = 4
$ws.Range("B31").Value = 11
$ws.Range("B34").Value = 10
$ws.Range("B37").Value = 11
$ws.Range("B38").Value = 26
$ws.Range("B40").Value = 19
$ws.Range("B42").Value = 44
$ws.Range("B45").Value = 101
$ws.Range("B46").Value = 81
$ws.Range("B47").Value = 119
$ws.Range("B48").Value = 118
$ws.Range("B49").Value = 188
$ws.Range("B50").Value = 162
$ws.Range("B51").Value = 385
$ws.Range("B52").Value = 370
$ws.Range("B53").Value = 426
$ws.Range("B54").Value = 471
$ws.Range("B55").Value = 553
$ws.Range("B56").Value = 735
$ws.Range("B57").Value = 883
$ws.Range("B58").Value = 1397
$ws.Range("B59").Value = 2355
$ws.Range("B60").Value = 2538
$ws.Range("B61").Value = 3358
$ws.Range("B62").Value = 5254
$ws.Range("B63").Value = 6606
$ws.Range("B64").Value = 7438
$ws.Range("B65").Value = 9100
$ws.Range("B66").Value = 10742
$ws.Range("B67").Value = 8945
$ws.Range("B68").Value = 8790
$ws.Range("B69").Value = 10483
$ws.Range("B70").Value = 12509
$ws.Range("B71").Value = 12365
$ws.Range("B73").Value = 14306
$ws.Range("B74").Value = 10409
$ws.Range("B75").Value = 10661
$ws.Range("B76").Value = 15546
$ws.Range("B77").Value = 14804
$ws.Range("B78").Value = 14256
$ws.Range("B79").Value = 14305
$ws.Range("B80").Value = 14976
$ws.Range("B81").Value = 11096
$ws.Range("B82").Value = 10479
$ws.Range("B83").Value = 15666
$ws.Range("B84").Value = 15052
$ws.Range("B85").Value = 13384
$ws.Range("B86").Value = 11014
$ws.Range("B87").Value = 9952
$ws.Range("B88").Value = 5271
$ws.Range("B89").Value = 2894
$ws.Range("B90").Value = 1606
$ws.Range("B91").Value = 386
$ws.Range("B92").Value = 126
$ws.Range("B93").Value = 56
$ws.Range("B94").Value = 13

# --- A84 loses its "latest-day" highlight now that later rows take over ---
# Copy the (unhighlighted) format of A83 onto A84, then restore A84's own date.
$ws.Range("A83").Copy($ws.Range("A84"))
$excel.CutCopyMode = $false
$ws.Range("A84").Value = 43924

# --- Append the new data row (row 95) with the highlight that A84 used to carry ---
$ws.Range("A94").Copy($ws.Range("A95"))
$excel.CutCopyMode = $false
$ws.Range("A95").Value = 43935
$ws.Range("B95").Value = 2

# --- View state: scroll the grid so the new tail rows are visible, cursor on A84 ---
$excel.ActiveWindow.ScrollRow = 78
$ws.Range("A84").Select() | Out-Null
